$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1640
$ws.Range("F4").Value = 9438
$ws.Range("F5").Value = 724
$ws.Range("F7").Value = 201
$ws.Range("F8").Value = 319
$ws.Range("F10").Value = 69
$ws.Range("F11").Value = 1615
$ws.Range("F12").Value = 1389
$ws.Range("F15").Value = 1457
$ws.Range("F17").Value = 296
$ws.Range("F19").Value = 132
$ws.Range("F20").Value = 81
$ws.Range("F21").Value = 367
$ws.Range("F22").Value = 1103
$ws.Range("F26").Value = 272
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = 251
$ws.Range("F29").Value = 73
$ws.Range("F30").Value = 603
$ws.Range("F32").Value = 2
$ws.Range("F33").Value = 159
$ws.Range("F36").Value = 6
$ws.Range("F37").Value = 2
$ws.Range("F38").Value = 218
$ws.Range("F39").Value = 594
$ws.Range("F40").Value = 2
$ws.Range("F42").Value = 731
$ws.Range("F43").Value = 4
$ws.Range("F44").Value = 270
$ws.Range("F45").Value = 6
$ws.Range("F46").Value = 59
$ws.Range("F47").Value = 4

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 15
$ws.Range("F6").Value = 64
$ws.Range("F15").Value = 44
$ws.Range("F17").Value = 126
$ws.Range("F18").Value = 965
$ws.Range("F20").Value = 1058
$ws.Range("F21").Value = 276
$ws.Range("F23").Value = 5
$ws.Range("F25").Value = 300
$ws.Range("F31").Value = 167
$ws.Range("F35").Value = 111
$ws.Range("F37").Value = 8

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 767
$ws.Range("F5").Value = 359
$ws.Range("F7").Value = 2342
$ws.Range("F8").Value = 3531
$ws.Range("F9").Value = 2
$ws.Range("F11").Value = 43
$ws.Range("F12").Value = 82

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1640
$ws.Range("F5").Value = 9438
$ws.Range("F7").Value = 3531
$ws.Range("F8").Value = 724
$ws.Range("F11").Value = 319
$ws.Range("F13").Value = 69
$ws.Range("F14").Value = 1389
$ws.Range("F16").Value = 82
$ws.Range("F19").Value = 296
$ws.Range("F20").Value = 132
$ws.Range("F21").Value = 367
$ws.Range("F22").Value = 1103
$ws.Range("F26").Value = 44
$ws.Range("F29").Value = 272
$ws.Range("F31").Value = 1058
$ws.Range("F32").Value = 276
$ws.Range("F33").Value = 73
$ws.Range("F34").Value = 603
$ws.Range("F37").Value = 159
$ws.Range("F39").Value = 300
$ws.Range("F40").Value = 218
$ws.Range("F42").Value = 594
$ws.Range("F43").Value = 731
$ws.Range("F45").Value = 167
$ws.Range("F46").Value = 270
$ws.Range("F47").Value = 111
$ws.Range("F48").Value = 59
